$d = $word.ActiveDocument

# Locate the paragraph that holds the "//margin-bottom: 0;" run so the new
# material can be anchored right after it.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*//margin-bottom: 0;*") {
        $targetIndex = $i
    }
}

# Insert two new blank paragraphs after it (re-fetching the paragraph by
# index each time so the Range used reflects the just-inserted paragraph
# mark rather than staying anchored to the original, already-collapsed
# range).
$r = $d.Paragraphs.Item($targetIndex).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIndex + 1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIndex + 2).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# The paragraph that was just created is now empty; fill it with the new
# text, leaving the original trailing empty paragraph untouched after it.
$newTextPara = $d.Paragraphs.Item($targetIndex + 3)
$newTextPara.Range.InsertAfter('data-html="true" => can be styled using html (css?)')
